# dataKichua.xlsx edit: add a new "Eres/ahora" level row and tweak the
# existing shared strings + column widths to match the updated puzzle data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing word-list strings in row 2 (B2/C2) to the new
# "word - hint" format.
$ws.Range("B2").Value2 = "Eres - heres,tanta,gente,quién,ahora - ajora"
$ws.Range("C2").Value2 = "Eres - heres,ahora - ajora"

# Append a new data row (row 3) that reuses the song name (A) and the
# freshly updated word lists (B/C), with the next level number in D.
$ws.Range("A3").Value2 = "alejandro sanz no es lo mismo"
$ws.Range("B3").Value2 = "Eres - heres,tanta,gente,quién,ahora - ajora"
$ws.Range("C3").Value2 = "Eres - heres,ahora - ajora"
$ws.Range("D3").Value2 = 2

# Widen columns B and C to fit the longer text.
$ws.Columns.Item(2).ColumnWidth = 38
$ws.Columns.Item(3).ColumnWidth = 31

# Move the active selection to D2, matching the saved view state.
$ws.Range("D2").Select()
